# AdventStats.xlsx - "2023 day 11 done"
#
# - "2023" sheet: refresh the day-by-day leaderboard numbers (days 1-10,
#   rows 2-11) with the latest "My 1"/"My 2" counts, and fill in day 11
#   (row 12), which had previously been left blank.
# - "Overall" sheet: mark 2023 day 11 part 1 (columns AP:AS, row 12) as
#   succeeded ("s") - it was previously "t" (not done yet).
# - Selection moves on to the next row/day on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("2023")
$ws2 = $wb.Worksheets.Item("Overall")

# --- "2023" sheet: updated totals for days already in progress ---
$ws1.Range("B2").Value  = 200872
$ws1.Range("C2").Value  = 63434

$ws1.Range("B3").Value  = 168562
$ws1.Range("C3").Value  = 7606

$ws1.Range("B4").Value  = 110983
$ws1.Range("C4").Value  = 16491

$ws1.Range("B5").Value  = 109681
$ws1.Range("C5").Value  = 14455

$ws1.Range("B6").Value  = 66368
$ws1.Range("C6").Value  = 25819

$ws1.Range("B7").Value  = 84067
$ws1.Range("C7").Value  = 1354

$ws1.Range("B8").Value  = 65148
$ws1.Range("C8").Value  = 5769

$ws1.Range("B9").Value  = 56966
$ws1.Range("C9").Value  = 11567

$ws1.Range("B10").Value = 55206
$ws1.Range("C10").Value = 746

$ws1.Range("B11").Value = 28621
$ws1.Range("C11").Value = 14324

# --- day 11 (row 12) just finished: fill in the previously empty cells ---
# Row 12 formulas (D12, H12) gate on ISBLANK(B12) / ISBLANK(C12); set the
# other inputs first and touch B12/C12 last so every dependent formula in
# the row picks up the new values on recalculation.
$ws1.Range("E12").Value = 12710
$ws1.Range("F12").Value = 11612
$ws1.Range("C12").Value = 1978
$ws1.Range("B12").Value = 12059
$ws1.Range("C12").Value = 1978

# --- "Overall" sheet: 2023 (row 12) day 11 (cols AP:AS) part 1 succeeded ---
$ws2.Range("AP12").Value = "s"
$ws2.Range("AQ12").Value = "s"
$ws2.Range("AR12").Value = "s"
$ws2.Range("AS12").Value = "s"

# --- selection state: move to the next cell on each sheet ---
$ws2.Range("AT12").Select()

$ws1.Activate()
$ws1.Range("B13").Select()
